$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "30.334.90"
$c.Style = $s
$ws.Range("E2").Value = "  -2.64%  "

$c = $ws.Range("D3")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.941.31"
$c.Style = $s
$ws.Range("E3").Value = "  -2.42%  "

$c = $ws.Range("D4")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = $s
$ws.Range("E4").Value = "  +0.28%  "

$c = $ws.Range("D5")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "251.25"
$c.Style = $s
$ws.Range("E5").Value = "  -1.52%  "

$c = $ws.Range("D6")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.7187"
$c.Style = $s
$ws.Range("E6").Value = "  -9.26%  "

$c = $ws.Range("D7")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = $s
$ws.Range("E7").Value = "  +0.23%  "

$c = $ws.Range("D8")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.3331"
$c.Style = $s
$ws.Range("E8").Value = "  -4.76%  "

$c = $ws.Range("D9")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "28.66"
$c.Style = $s
$ws.Range("E9").Value = "  +2.16%  "

$c = $ws.Range("D10")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.07369"
$c.Style = $s
$ws.Range("E10").Value = "  +5.34%  "

$c = $ws.Range("D11")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.8143"
$c.Style = $s
$ws.Range("E11").Value = "  -3.74%  "

$c = $ws.Range("D12")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.08135"
$c.Style = $s
$ws.Range("E12").Value = "  -0.67%  "

$c = $ws.Range("D13")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.941.49"
$c.Style = $s
$ws.Range("E13").Value = "  -2.43%  "

$c = $ws.Range("D14")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.481"
$c.Style = $s
$ws.Range("E14").Value = "  -2.36%  "

$c = $ws.Range("D15")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "95.03"
$c.Style = $s
$ws.Range("E15").Value = "  -5.40%  "

$c = $ws.Range("D16")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "14.89"
$c.Style = $s
$ws.Range("E16").Value = "  -3.28%  "

$c = $ws.Range("D17")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.000008438"
$c.Style = $s
$ws.Range("E17").Value = "  +6.46%  "

$c = $ws.Range("D18")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "30.369.60"
$c.Style = $s
$ws.Range("E18").Value = "  -2.55%  "

$c = $ws.Range("D19")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "253.35"
$c.Style = $s
$ws.Range("E19").Value = "  -7.32%  "

$c = $ws.Range("D20")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.883"
$c.Style = $s
$ws.Range("E20").Value = "  +0.20%  "

$c = $ws.Range("D21")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.196.74"
$c.Style = $s
$ws.Range("E21").Value = "  -2.24%  "

$c = $ws.Range("D22")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = $s
$ws.Range("E22").Value = "  +0.10%  "

$ws.Range("E23").Value = "  +0.11%  "

$c = $ws.Range("D24")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.977"
$c.Style = $s
$ws.Range("E24").Value = "  -1.08%  "

$c = $ws.Range("D25")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "9.843"
$c.Style = $s
$ws.Range("E25").Value = "  -1.39%  "

$c = $ws.Range("D26")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "162.86"
$c.Style = $s
$ws.Range("E26").Value = "  -1.97%  "

$c = $ws.Range("D27")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.411"
$c.Style = $s
$ws.Range("E27").Value = "  +3.12%  "

$c = $ws.Range("D28")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "19.36"
$c.Style = $s
$ws.Range("E28").Value = "  -2.79%  "

$c = $ws.Range("D29")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.1318"
$c.Style = $s
$ws.Range("E29").Value = "  -12.24%  "

$c = $ws.Range("D30")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.570"
$c.Style = $s
$ws.Range("E30").Value = "  -1.75%  "

$ws.Range("E31").Value = "  -0.56%  "

$c = $ws.Range("D32")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "4.453"
$c.Style = $s
$ws.Range("E32").Value = "  -2.86%  "

$c = $ws.Range("D33")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "4.243"
$c.Style = $s
$ws.Range("E33").Value = "  -3.64%  "

$c = $ws.Range("D34")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.05273"
$c.Style = $s
$ws.Range("E34").Value = "  +0.71%  "

$c = $ws.Range("D35")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.299"
$c.Style = $s
$ws.Range("E35").Value = "  +6.42%  "

$c = $ws.Range("D36")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.7558"
$c.Style = $s
$ws.Range("E36").Value = "  -2.96%  "

$c = $ws.Range("D37")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.743"
$c.Style = $s
$ws.Range("E37").Value = "  -0.90%  "

$c = $ws.Range("D38")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.01991"
$c.Style = $s
$ws.Range("E38").Value = "  -0.44%  "

$c = $ws.Range("D39")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.851"
$c.Style = $s
$ws.Range("E39").Value = "  -0.99%  "

$c = $ws.Range("D40")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "81.13"
$c.Style = $s
$ws.Range("E40").Value = "  +1.86%  "

$c = $ws.Range("D41")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.609"
$c.Style = $s
$ws.Range("E41").Value = "  -0.30%  "

$c = $ws.Range("D42")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.4557"
$c.Style = $s
$ws.Range("E42").Value = "  -2.46%  "

$c = $ws.Range("D43")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.038"
$c.Style = $s
$ws.Range("E43").Value = "  -3.80%  "

$c = $ws.Range("D44")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.8479"
$c.Style = $s
$ws.Range("E44").Value = "  -0.18%  "

$ws.Range("E45").Value = "  +0.21%  "

$c = $ws.Range("D46")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "102.82"
$c.Style = $s
$ws.Range("E46").Value = "  -1.82%  "

$c = $ws.Range("D47")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "9.877"
$c.Style = $s
$ws.Range("E47").Value = "  +0.25%  "

$c = $ws.Range("D48")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.498"
$c.Style = $s
$ws.Range("E48").Value = "  -2.36%  "

$c = $ws.Range("D49")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "36.84"
$c.Style = $s
$ws.Range("E49").Value = "  +0.33%  "

$c = $ws.Range("D50")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.4191"
$c.Style = $s
$ws.Range("E50").Value = "  -2.35%  "

$c = $ws.Range("D51")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.505"
$c.Style = $s
$ws.Range("E51").Value = "  -2.01%  "
